# Update the "想去人数" (want-to-go count) figures for the generated
# gh-pages data dump. The same set of events is listed on both the
# "展览" (Exhibitions) sheet and the "全部类型" (All types) sheet, so the
# corresponding cells need to be updated on both sheets.

$wb = $excel.ActiveWorkbook

$updates = @{
    "展览"     = @{ "F2" = 14895; "F3" = 18597; "F5" = 115; "F13" = 51; "F17" = 1421; "F20" = 86; "F21" = 231; "F22" = 7715; "F26" = 1224; "F28" = 5972; "F33" = 261; "F34" = 5329 }
    "全部类型" = @{ "F2" = 14895; "F3" = 18597; "F5" = 115; "F13" = 51; "F17" = 1421; "F21" = 86; "F22" = 231; "F23" = 7715; "F27" = 1224; "F31" = 5972; "F36" = 261; "F37" = 5329 }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $cells = $updates[$sheetName]
    foreach ($addr in $cells.Keys) {
        $ws.Range($addr).Value = $cells[$addr]
    }
}
